$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A103").Value = 46034.76892335648
$ws.Range("A103").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B103").Value = "MA001BOR"
$ws.Range("C103").Value = "Borno"
$ws.Range("D103").Value = "Lashe Money"
$ws.Range("E103").Value = "Soya Beans"
$ws.Range("F103").Value = 56000
$ws.Range("F103").NumberFormat = "0.00"
$ws.Range("G103").Value = 115
$ws.Range("G103").Style = "Normal"
$ws.Range("H103").Value = 486.95652173913
$ws.Range("H103").Style = "Normal"
$ws.Range("I103").Value = "high"
$ws.Range("J103").Value = "New"
$ws.Range("K103").Value = 489.95652173913
$ws.Range("K103").Style = "Normal"

$ws.Range("A104").Value = 46034.76052487268
$ws.Range("A104").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B104").Value = "MA001BOR"
$ws.Range("C104").Value = "Borno"
$ws.Range("D104").Value = "Lashe Money"
$ws.Range("E104").Value = "Honeybeans"
$ws.Range("F104").Value = 62000
$ws.Range("F104").NumberFormat = "0.00"
$ws.Range("G104").Value = 103
$ws.Range("G104").Style = "Normal"
$ws.Range("H104").Value = 601.941747572815
$ws.Range("H104").Style = "Normal"
$ws.Range("I104").Value = "high"
$ws.Range("J104").Value = "New"
$ws.Range("K104").Value = 604.941747572815
$ws.Range("K104").Style = "Normal"

$ws.Range("A105").Value = 46034.75988082176
$ws.Range("A105").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B105").Value = "MA001BOR"
$ws.Range("C105").Value = "Borno"
$ws.Range("D105").Value = "Lashe Money"
$ws.Range("E105").Value = "Cowpea Brown"
$ws.Range("F105").Value = 58000
$ws.Range("F105").NumberFormat = "0.00"
$ws.Range("G105").Value = 105
$ws.Range("G105").Style = "Normal"
$ws.Range("H105").Value = 552.380952380952
$ws.Range("H105").Style = "Normal"
$ws.Range("I105").Value = "high"
$ws.Range("J105").Value = "New"
$ws.Range("K105").Value = 555.380952380952
$ws.Range("K105").Style = "Normal"

$ws.Range("A106").Value = 46034.75926842593
$ws.Range("A106").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B106").Value = "MA001BOR"
$ws.Range("C106").Value = "Borno"
$ws.Range("D106").Value = "Lashe Money"
$ws.Range("E106").Value = "Cowpea White"
$ws.Range("F106").Value = 61000
$ws.Range("F106").NumberFormat = "0.00"
$ws.Range("G106").Value = 105
$ws.Range("G106").Style = "Normal"
$ws.Range("H106").Value = 580.952380952381
$ws.Range("H106").Style = "Normal"
$ws.Range("I106").Value = "high"
$ws.Range("J106").Value = "New"
$ws.Range("K106").Value = 583.952380952381
$ws.Range("K106").Style = "Normal"

$ws.Range("A107").Value = 46034.7585496875
$ws.Range("A107").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B107").Value = "MA001BOR"
$ws.Range("C107").Value = "Borno"
$ws.Range("D107").Value = "Lashe Money"
$ws.Range("E107").Value = "Rice Paddy"
$ws.Range("F107").Value = 35000
$ws.Range("F107").NumberFormat = "0.00"
$ws.Range("G107").Value = 95
$ws.Range("G107").Style = "Normal"
$ws.Range("H107").Value = 368.421052631578
$ws.Range("H107").Style = "Normal"
$ws.Range("I107").Value = "high"
$ws.Range("J107").Value = "New"
$ws.Range("K107").Value = 371.421052631578
$ws.Range("K107").Style = "Normal"

$ws.Range("A108").Value = 46037.76779003472
$ws.Range("A108").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B108").Value = "MH001GMB"
$ws.Range("C108").Value = "Gombe"
$ws.Range("D108").Value = "Kumo"
$ws.Range("E108").Value = "Soya Beans"
$ws.Range("F108").Value = 52000
$ws.Range("F108").NumberFormat = "0.00"
$ws.Range("G108").Value = 100
$ws.Range("G108").Style = "Normal"
$ws.Range("H108").Value = 520
$ws.Range("H108").Style = "Normal"
$ws.Range("I108").Value = "medium"
$ws.Range("J108").Value = "New"
$ws.Range("K108").Value = 523
$ws.Range("K108").Style = "Normal"

$ws.Range("A109").Value = 46037.76694519676
$ws.Range("A109").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B109").Value = "MH001GMB"
$ws.Range("C109").Value = "Gombe"
$ws.Range("D109").Value = "Kumo"
$ws.Range("E109").Value = "Sorghum Red"
$ws.Range("F109").Value = 22000
$ws.Range("F109").NumberFormat = "0.00"
$ws.Range("G109").Value = 100
$ws.Range("G109").Style = "Normal"
$ws.Range("H109").Value = 220
$ws.Range("H109").Style = "Normal"
$ws.Range("I109").Value = "low"
$ws.Range("J109").Value = "New"
$ws.Range("K109").Value = 223
$ws.Range("K109").Style = "Normal"

$ws.Range("A110").Value = 46037.76623640046
$ws.Range("A110").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B110").Value = "MH001GMB"
$ws.Range("C110").Value = "Gombe"
$ws.Range("D110").Value = "Kumo"
$ws.Range("E110").Value = "Rice Paddy"
$ws.Range("F110").Value = 27000
$ws.Range("F110").NumberFormat = "0.00"
$ws.Range("G110").Value = 70
$ws.Range("G110").Style = "Normal"
$ws.Range("H110").Value = 385.714285714285
$ws.Range("H110").Style = "Normal"
$ws.Range("I110").Value = "medium"
$ws.Range("J110").Value = "New"
$ws.Range("K110").Value = 388.714285714285
$ws.Range("K110").Style = "Normal"

$ws.Range("A111").Value = 46037.76526109954
$ws.Range("A111").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B111").Value = "MH001GMB"
$ws.Range("C111").Value = "Gombe"
$ws.Range("D111").Value = "Kumo"
$ws.Range("E111").Value = "Millet"
$ws.Range("F111").Value = 22000
$ws.Range("F111").NumberFormat = "0.00"
$ws.Range("G111").Value = 100
$ws.Range("G111").Style = "Normal"
$ws.Range("H111").Value = 220
$ws.Range("H111").Style = "Normal"
$ws.Range("I111").Value = "medium"
$ws.Range("J111").Value = "New"
$ws.Range("K111").Value = 223
$ws.Range("K111").Style = "Normal"

$ws.Range("A112").Value = 46037.76440118055
$ws.Range("A112").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B112").Value = "MH001GMB"
$ws.Range("C112").Value = "Gombe"
$ws.Range("D112").Value = "Kumo"
$ws.Range("E112").Value = "Maize White"
$ws.Range("F112").Value = 22000
$ws.Range("F112").NumberFormat = "0.00"
$ws.Range("G112").Value = 95
$ws.Range("G112").Style = "Normal"
$ws.Range("H112").Value = 231.578947368421
$ws.Range("H112").Style = "Normal"
$ws.Range("I112").Value = "medium"
$ws.Range("J112").Value = "New"
$ws.Range("K112").Value = 234.578947368421
$ws.Range("K112").Style = "Normal"

$ws.Range("A113").Value = 46037.76353651621
$ws.Range("A113").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B113").Value = "MH001GMB"
$ws.Range("C113").Value = "Gombe"
$ws.Range("D113").Value = "Kumo"
$ws.Range("E113").Value = "Groundnut Gargaja"
$ws.Range("F113").Value = 90000
$ws.Range("F113").NumberFormat = "0.00"
$ws.Range("G113").Value = 85
$ws.Range("G113").Style = "Normal"
$ws.Range("H113").Value = 1058.82352941176
$ws.Range("H113").Style = "Normal"
$ws.Range("I113").Value = "medium"
$ws.Range("J113").Value = "New"
$ws.Range("K113").Value = 1061.82352941176
$ws.Range("K113").Style = "Normal"

$ws.Range("A114").Value = 46037.76272434027
$ws.Range("A114").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B114").Value = "MH001GMB"
$ws.Range("C114").Value = "Gombe"
$ws.Range("D114").Value = "Kumo"
$ws.Range("E114").Value = "Cowpea White"
$ws.Range("F114").Value = 51000
$ws.Range("F114").NumberFormat = "0.00"
$ws.Range("G114").Value = 80
$ws.Range("G114").Style = "Normal"
$ws.Range("H114").Value = 637.5
$ws.Range("H114").Style = "Normal"
$ws.Range("I114").Value = "medium"
$ws.Range("J114").Value = "New"
$ws.Range("K114").Value = 640.5
$ws.Range("K114").Style = "Normal"

$ws.Range("A115").Value = 46038.62758684028
$ws.Range("A115").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B115").Value = "IS001KDN"
$ws.Range("C115").Value = "kaduna"
$ws.Range("D115").Value = "giwa"
$ws.Range("E115").Value = "Maize"
$ws.Range("F115").Value = 23000
$ws.Range("F115").Style = "Normal"
$ws.Range("G115").Value = 100
$ws.Range("G115").Style = "Normal"
$ws.Range("H115").Value = 230
$ws.Range("H115").Style = "Normal"
$ws.Range("I115").Value = "high"
$ws.Range("J115").Value = "New"
$ws.Range("K115").Value = 233
$ws.Range("K115").Style = "Normal"

$ws.Range("A116").Value = 46038.62829288194
$ws.Range("A116").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B116").Value = "IS001KDN"
$ws.Range("C116").Value = "kaduna"
$ws.Range("D116").Value = "giwa"
$ws.Range("E116").Value = "Sorghum"
$ws.Range("F116").Value = 24000
$ws.Range("F116").Style = "Normal"
$ws.Range("G116").Value = 100
$ws.Range("G116").Style = "Normal"
$ws.Range("H116").Value = 240
$ws.Range("H116").Style = "Normal"
$ws.Range("I116").Value = "high"
$ws.Range("J116").Value = "New"
$ws.Range("K116").Value = 243
$ws.Range("K116").Style = "Normal"

$ws.Range("A117").Value = 46038.62888636574
$ws.Range("A117").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B117").Value = "IS001KDN"
$ws.Range("C117").Value = "kaduna"
$ws.Range("D117").Value = "giwa"
$ws.Range("E117").Value = "Sorghum White"
$ws.Range("F117").Value = 25000
$ws.Range("F117").Style = "Normal"
$ws.Range("G117").Value = 100
$ws.Range("G117").Style = "Normal"
$ws.Range("H117").Value = 250
$ws.Range("H117").Style = "Normal"
$ws.Range("I117").Value = "high"
$ws.Range("J117").Value = "New"
$ws.Range("K117").Value = 253
$ws.Range("K117").Style = "Normal"

$ws.Range("A118").Value = 46038.62943975694
$ws.Range("A118").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B118").Value = "IS001KDN"
$ws.Range("C118").Value = "kaduna"
$ws.Range("D118").Value = "giwa"
$ws.Range("E118").Value = "Soya Beans"
$ws.Range("F118").Value = 50000
$ws.Range("F118").Style = "Normal"
$ws.Range("G118").Value = 100
$ws.Range("G118").Style = "Normal"
$ws.Range("H118").Value = 500
$ws.Range("H118").Style = "Normal"
$ws.Range("I118").Value = "high"
$ws.Range("J118").Value = "New"
$ws.Range("K118").Value = 503
$ws.Range("K118").Style = "Normal"

$ws.Range("A119").Value = 46038.69045054398
$ws.Range("A119").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B119").Value = "IS001KDN"
$ws.Range("C119").Value = "kaduna"
$ws.Range("D119").Value = "pambegua"
$ws.Range("E119").Value = "Rice Paddy"
$ws.Range("F119").Value = 34000
$ws.Range("F119").Style = "Normal"
$ws.Range("G119").Value = 350
$ws.Range("G119").Style = "Normal"
$ws.Range("H119").Value = 97.1428571428571
$ws.Range("H119").Style = "Normal"
$ws.Range("I119").Value = "medium"
$ws.Range("J119").Value = "New"
$ws.Range("K119").Value = 100.142857142857
$ws.Range("K119").Style = "Normal"

$ws.Range("A120").Value = 46038.6971699537
$ws.Range("A120").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B120").Value = "IS001KDN"
$ws.Range("C120").Value = "kaduna"
$ws.Range("D120").Value = "pambegua"
$ws.Range("E120").Value = "Maize White"
$ws.Range("F120").Value = 22500
$ws.Range("F120").Style = "Normal"
$ws.Range("G120").Value = 235
$ws.Range("G120").Style = "Normal"
$ws.Range("H120").Value = 95.7446808510638
$ws.Range("H120").Style = "Normal"
$ws.Range("I120").Value = "high"
$ws.Range("J120").Value = "New"
$ws.Range("K120").Value = 98.7446808510638
$ws.Range("K120").Style = "Normal"

$ws.Range("A121").Value = 46038.70371695602
$ws.Range("A121").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B121").Value = "IS001KDN"
$ws.Range("C121").Value = "kaduna"
$ws.Range("D121").Value = "pambegua"
$ws.Range("E121").Value = "Soya Beans"
$ws.Range("F121").Value = 51000
$ws.Range("F121").Style = "Normal"
$ws.Range("G121").Value = 530
$ws.Range("G121").Style = "Normal"
$ws.Range("H121").Value = 96.2264150943396
$ws.Range("H121").Style = "Normal"
$ws.Range("I121").Value = "medium"
$ws.Range("J121").Value = "Old"
$ws.Range("K121").Value = 99.2264150943396
$ws.Range("K121").Style = "Normal"

$ws.Range("A122").Value = 46038.70497664352
$ws.Range("A122").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B122").Value = "IS001KDN"
$ws.Range("C122").Value = "kaduna"
$ws.Range("D122").Value = "pambegua"
$ws.Range("E122").Value = "Sorghum White"
$ws.Range("F122").Value = 25000
$ws.Range("F122").Style = "Normal"
$ws.Range("G122").Value = 275
$ws.Range("G122").Style = "Normal"
$ws.Range("H122").Value = 90.9090909090909
$ws.Range("H122").Style = "Normal"
$ws.Range("I122").Value = "medium"
$ws.Range("J122").Value = "Old"
$ws.Range("K122").Value = 93.9090909090909
$ws.Range("K122").Style = "Normal"

$ws.Range("A123").Value = 46039.78372648148
$ws.Range("A123").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B123").Value = "MH001GMB"
$ws.Range("C123").Value = "gombe"
$ws.Range("D123").Value = "biliri"
$ws.Range("E123").Value = "Cowpea White"
$ws.Range("F123").Value = 43000
$ws.Range("F123").Style = "Normal"
$ws.Range("G123").Value = 80
$ws.Range("G123").Style = "Normal"
$ws.Range("H123").Value = 537.5
$ws.Range("H123").Style = "Normal"
$ws.Range("I123").Value = "medium"
$ws.Range("J123").Value = "New"
$ws.Range("K123").Value = 540.5
$ws.Range("K123").Style = "Normal"

$ws.Range("A124").Value = 46039.78439196759
$ws.Range("A124").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B124").Value = "MH001GMB"
$ws.Range("C124").Value = "gombe"
$ws.Range("D124").Value = "biliri"
$ws.Range("E124").Value = "Groundnut Gargaja"
$ws.Range("F124").Value = 85000
$ws.Range("F124").Style = "Normal"
$ws.Range("G124").Value = 85
$ws.Range("G124").Style = "Normal"
$ws.Range("H124").Value = 1000
$ws.Range("H124").Style = "Normal"
$ws.Range("I124").Value = "medium"
$ws.Range("J124").Value = "New"
$ws.Range("K124").Value = 1003
$ws.Range("K124").Style = "Normal"

$ws.Range("A125").Value = 46039.78497951389
$ws.Range("A125").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B125").Value = "MH001GMB"
$ws.Range("C125").Value = "gombe"
$ws.Range("D125").Value = "biliri"
$ws.Range("E125").Value = "Groundut Kampala"
$ws.Range("F125").Value = 80000
$ws.Range("F125").Style = "Normal"
$ws.Range("G125").Value = 85
$ws.Range("G125").Style = "Normal"
$ws.Range("H125").Value = 941.176470588235
$ws.Range("H125").Style = "Normal"
$ws.Range("I125").Value = "low"
$ws.Range("J125").Value = "New"
$ws.Range("K125").Value = 944.176470588235
$ws.Range("K125").Style = "Normal"

$ws.Range("A126").Value = 46039.78557793982
$ws.Range("A126").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B126").Value = "MH001GMB"
$ws.Range("C126").Value = "gombe"
$ws.Range("D126").Value = "biliri"
$ws.Range("E126").Value = "Maize White"
$ws.Range("F126").Value = 22000
$ws.Range("F126").Style = "Normal"
$ws.Range("G126").Value = 95
$ws.Range("G126").Style = "Normal"
$ws.Range("H126").Value = 231.578947368421
$ws.Range("H126").Style = "Normal"
$ws.Range("I126").Value = "medium"
$ws.Range("J126").Value = "New"
$ws.Range("K126").Value = 234.578947368421
$ws.Range("K126").Style = "Normal"

$ws.Range("A127").Value = 46039.78628592593
$ws.Range("A127").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B127").Value = "MH001GMB"
$ws.Range("C127").Value = "gombe"
$ws.Range("D127").Value = "biliri"
$ws.Range("E127").Value = "Millet"
$ws.Range("F127").Value = 25000
$ws.Range("F127").Style = "Normal"
$ws.Range("G127").Value = 100
$ws.Range("G127").Style = "Normal"
$ws.Range("H127").Value = 250
$ws.Range("H127").Style = "Normal"
$ws.Range("I127").Value = "low"
$ws.Range("J127").Value = "New"
$ws.Range("K127").Value = 253
$ws.Range("K127").Style = "Normal"

$ws.Range("A128").Value = 46039.7871252662
$ws.Range("A128").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B128").Value = "MH001GMB"
$ws.Range("C128").Value = "gombe"
$ws.Range("D128").Value = "biliri"
$ws.Range("E128").Value = "Rice Paddy"
$ws.Range("F128").Value = 27000
$ws.Range("F128").Style = "Normal"
$ws.Range("G128").Value = 70
$ws.Range("G128").Style = "Normal"
$ws.Range("H128").Value = 385.714285714285
$ws.Range("H128").Style = "Normal"
$ws.Range("I128").Value = "medium"
$ws.Range("J128").Value = "New"
$ws.Range("K128").Value = 388.714285714285
$ws.Range("K128").Style = "Normal"

$ws.Range("A129").Value = 46039.78781471065
$ws.Range("A129").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B129").Value = "MH001GMB"
$ws.Range("C129").Value = "gombe"
$ws.Range("D129").Value = "biliri"
$ws.Range("E129").Value = "Sorghum Red"
$ws.Range("F129").Value = 28000
$ws.Range("F129").Style = "Normal"
$ws.Range("G129").Value = 100
$ws.Range("G129").Style = "Normal"
$ws.Range("H129").Value = 280
$ws.Range("H129").Style = "Normal"
$ws.Range("I129").Value = "low"
$ws.Range("J129").Value = "New"
$ws.Range("K129").Value = 283
$ws.Range("K129").Style = "Normal"

$ws.Range("A130").Value = 46039.78864159722
$ws.Range("A130").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B130").Value = "MH001GMB"
$ws.Range("C130").Value = "gombe"
$ws.Range("D130").Value = "biliri"
$ws.Range("E130").Value = "Soya Beans"
$ws.Range("F130").Value = 53000
$ws.Range("F130").Style = "Normal"
$ws.Range("G130").Value = 100
$ws.Range("G130").Style = "Normal"
$ws.Range("H130").Value = 530
$ws.Range("H130").Style = "Normal"
$ws.Range("I130").Value = "medium"
$ws.Range("J130").Value = "New"
$ws.Range("K130").Value = 533
$ws.Range("K130").Style = "Normal"

$ws.Range("A131").Value = 46040.85899197917
$ws.Range("A131").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B131").Value = "AU001YBE"
$ws.Range("C131").Value = "yobe"
$ws.Range("D131").Value = "potiskum"
$ws.Range("E131").Value = "Cowpea White"
$ws.Range("F131").Value = 49000
$ws.Range("F131").Style = "Normal"
$ws.Range("G131").Value = 100
$ws.Range("G131").Style = "Normal"
$ws.Range("H131").Value = 490
$ws.Range("H131").Style = "Normal"
$ws.Range("I131").Value = "high"
$ws.Range("J131").Value = "New"
$ws.Range("K131").Value = 493
$ws.Range("K131").Style = "Normal"

$ws.Range("A132").Value = 46040.85958349537
$ws.Range("A132").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B132").Value = "AU001YBE"
$ws.Range("C132").Value = "yobe"
$ws.Range("D132").Value = "potiskum"
$ws.Range("E132").Value = "Groundnut Gargaja"
$ws.Range("F132").Value = 107000
$ws.Range("F132").Style = "Normal"
$ws.Range("G132").Value = 98
$ws.Range("G132").Style = "Normal"
$ws.Range("H132").Value = 1091.83673469387
$ws.Range("H132").Style = "Normal"
$ws.Range("I132").Value = "high"
$ws.Range("J132").Value = "New"
$ws.Range("K132").Value = 1094.83673469387
$ws.Range("K132").Style = "Normal"

$ws.Range("A133").Value = 46040.86052976852
$ws.Range("A133").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B133").Value = "AU001YBE"
$ws.Range("C133").Value = "yobe"
$ws.Range("D133").Value = "potiskum"
$ws.Range("E133").Value = "Maize White"
$ws.Range("F133").Value = 27000
$ws.Range("F133").Style = "Normal"
$ws.Range("G133").Value = 100
$ws.Range("G133").Style = "Normal"
$ws.Range("H133").Value = 270
$ws.Range("H133").Style = "Normal"
$ws.Range("I133").Value = "high"
$ws.Range("J133").Value = "New"
$ws.Range("K133").Value = 273
$ws.Range("K133").Style = "Normal"

$ws.Range("A134").Value = 46040.86148290509
$ws.Range("A134").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B134").Value = "AU001YBE"
$ws.Range("C134").Value = "yobe"
$ws.Range("D134").Value = "potiskum"
$ws.Range("E134").Value = "Millet"
$ws.Range("F134").Value = 30000
$ws.Range("F134").Style = "Normal"
$ws.Range("G134").Value = 105
$ws.Range("G134").Style = "Normal"
$ws.Range("H134").Value = 285.714285714285
$ws.Range("H134").Style = "Normal"
$ws.Range("I134").Value = "high"
$ws.Range("J134").Value = "New"
$ws.Range("K134").Value = 288.714285714285
$ws.Range("K134").Style = "Normal"

$ws.Range("A135").Value = 46040.86315017361
$ws.Range("A135").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B135").Value = "AU001YBE"
$ws.Range("C135").Value = "yobe"
$ws.Range("D135").Value = "potiskum"
$ws.Range("E135").Value = "Rice Processed"
$ws.Range("F135").Value = 900000
$ws.Range("F135").Style = "Normal"
$ws.Range("G135").Value = 115
$ws.Range("G135").Style = "Normal"
$ws.Range("H135").Value = 7826.08695652173
$ws.Range("H135").Style = "Normal"
$ws.Range("I135").Value = "high"
$ws.Range("J135").Value = "New"
$ws.Range("K135").Value = 7829.08695652173
$ws.Range("K135").Style = "Normal"

$ws.Range("A136").Value = 46040.86404246528
$ws.Range("A136").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B136").Value = "AU001YBE"
$ws.Range("C136").Value = "yobe"
$ws.Range("D136").Value = "potiskum"
$ws.Range("E136").Value = "Sorghum"
$ws.Range("F136").Value = 20000
$ws.Range("F136").Style = "Normal"
$ws.Range("G136").Value = 95
$ws.Range("G136").Style = "Normal"
$ws.Range("H136").Value = 210.526315789473
$ws.Range("H136").Style = "Normal"
$ws.Range("I136").Value = "high"
$ws.Range("J136").Value = "New"
$ws.Range("K136").Value = 213.526315789473
$ws.Range("K136").Style = "Normal"

$ws.Range("A137").Value = 46040.86605857639
$ws.Range("A137").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B137").Value = "AU001YBE"
$ws.Range("C137").Value = "yobe"
$ws.Range("D137").Value = "potiskum"
$ws.Range("E137").Value = "Sorghum Red"
$ws.Range("F137").Value = 23000
$ws.Range("F137").Style = "Normal"
$ws.Range("G137").Value = 95
$ws.Range("G137").Style = "Normal"
$ws.Range("H137").Value = 242.105263157894
$ws.Range("H137").Style = "Normal"
$ws.Range("I137").Value = "high"
$ws.Range("J137").Value = "New"
$ws.Range("K137").Value = 245.105263157894
$ws.Range("K137").Style = "Normal"

$ws.Range("A138").Value = 46040.868845
$ws.Range("A138").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B138").Value = "AU001YBE"
$ws.Range("C138").Value = "yobe"
$ws.Range("D138").Value = "potiskum"
$ws.Range("E138").Value = "Sorghum Yellow"
$ws.Range("F138").Value = 23000
$ws.Range("F138").Style = "Normal"
$ws.Range("G138").Value = 98
$ws.Range("G138").Style = "Normal"
$ws.Range("H138").Value = 234.69387755102
$ws.Range("H138").Style = "Normal"
$ws.Range("I138").Value = "high"
$ws.Range("J138").Value = "New"
$ws.Range("K138").Value = 237.69387755102
$ws.Range("K138").Style = "Normal"
